$wb = $excel.ActiveWorkbook

# --- Login sheet: adjust a few row heights ---
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Rows.Item(4).RowHeight = 19.5
$wsLogin.Rows.Item(12).RowHeight = 33
$wsLogin.Rows.Item(13).RowHeight = 33

# --- ChangePass sheet: add a new "Pass" column (old password test data) ---
$wsChangePass = $wb.Worksheets.Item("ChangePass")
$wsChangePass.Columns.Item(2).Insert()

# Copy formatting (style) for the new column from the last column's header cell,
# which already carries the plain "general" style used throughout column B.
for ($r = 1; $r -le 8; $r++) {
    $wsChangePass.Cells.Item(1, 10).Copy()
    $wsChangePass.Cells.Item($r, 2).PasteSpecial(-4122)
}

$wsChangePass.Range("B1").Value = "Pass"
$wsChangePass.Range("B2").Value = "adminmaster"
$wsChangePass.Range("B3").Value = "adminmaster"
$wsChangePass.Range("B4").Value = "adminmaster"
$wsChangePass.Range("B5").Value = "adminmaster"
$wsChangePass.Range("B6").Value = "adminmaster"
$wsChangePass.Range("B7").Value = "adminmaster"
$wsChangePass.Range("B8").Value = "adminmaster"

# ChangePass becomes the selected/active tab (moves away from SCartDel)
$wsChangePass.Activate()
